$wb = $excel.ActiveWorkbook

# Update the workbook window position (bookViews/workbookView)
$wb.Windows.Item(1).Left = 80
$wb.Windows.Item(1).Top = 460

# "All simulations" sheet: update simulation status values
$ws = $wb.Worksheets.Item("All simulations")
$ws.Range("J5").Value = "Running"
$ws.Range("J6").Value = "Done"

# Update the active selection on the sheet to reflect current work location
$ws.Range("J7").Select()

$wb.Save()
